$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the changed data values
$ws.Range("H2").Value = 5
$ws.Range("I2").Value = 197

# Update view: scroll back to A1 (remove topLeftCell="B1") and change selection
$ws.Range("A1").Select()
$ws.Range("I20").Select()
